$d = $word.ActiveDocument

$replacements = @(
    @{old="840×6=5040"; new="417×2=834"},
    @{old="434×2=868"; new="446×3=1338"},
    @{old="363×7=2541"; new="187×5=935"},
    @{old="298×8=2384"; new="288×2=576"},
    @{old="880×4=3520"; new="361×3=1083"},
    @{old="397×6=2382"; new="438×3=1314"},
    @{old="730×9=6570"; new="600×2=1200"},
    @{old="515×8=4120"; new="965×4=3860"},
    @{old="378×7=2646"; new="555×8=4440"},
    @{old="420×5=2100"; new="832×3=2496"},
    @{old="457×3=1371"; new="965×6=5790"},
    @{old="710×7=4970"; new="467×5=2335"},
    @{old="873×7=6111"; new="128×9=1152"},
    @{old="357×2=714"; new="302×5=1510"},
    @{old="380×8=3040"; new="885×7=6195"},
    @{old="799×6=4794"; new="200×9=1800"},
    @{old="721×7=5047"; new="322×6=1932"},
    @{old="242×2=484"; new="104×2=208"},
    @{old="423×3=1269"; new="777×7=5439"},
    @{old="494×6=2964"; new="216×6=1296"},
    @{old="526×7=3682"; new="757×4=3028"},
    @{old="726×7=5082"; new="651×9=5859"},
    @{old="254×5=1270"; new="733×5=3665"},
    @{old="431×6=2586"; new="565×4=2260"},
    @{old="978×6=5868"; new="858×8=6864"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
